$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '39.570.45'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').Value = '2.167.54'
$ws.Range('E3').Value = '  +3.54%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'229.61"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').Value = "'0.624"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('D7').Value = "'63.23"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.67%  '
$ws.Range('D9').Value = "'0.396"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.83%  '
$ws.Range('D10').Value = "'0.0860"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.58%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = "'16.12"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.74%  '
$ws.Range('D13').Value = '2.482.52'
$ws.Range('E13').Value = '  +3.36%  '
$ws.Range('D14').Value = "'22.34"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.58%  '
$ws.Range('E15').Value = '  +3.28%  '
$ws.Range('D16').Value = "'5.59"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.25%  '
$ws.Range('D17').Value = '2.156.09'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').Value = '39.576.31'
$ws.Range('E18').Value = '  +2.00%  '
$ws.Range('D19').Value = "'72.54"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.40%  '
$ws.Range('D20').Value = "'6.17"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('D21').Value = '0.0₃0855'
$ws.Range('E21').Value = '  +2.10%  '
$ws.Range('D22').Value = "'228.99"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').Value = "'2.35"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('D26').Value = "'9.69"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.96%  '
$ws.Range('D27').Value = "'173.01"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.21%  '
$ws.Range('D28').Value = "'0.139"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.40%  '
$ws.Range('D29').Value = "'1.43"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.01%  '
$ws.Range('D30').Value = "'19.72"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.70%  '
$ws.Range('E31').Value = '  +8.82%  '
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('D33').Value = "'4.66"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.54%  '
$ws.Range('D34').Value = "'4.83"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.66%  '
$ws.Range('D35').Value = "'7.11"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.81%  '
$ws.Range('E36').Value = '  +1.87%  '
$ws.Range('D37').Value = "'2.46"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.90%  '
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').Value = "'0.998"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.27%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0232"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.77%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = "'18.14"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('D42').Value = "'103.95"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.84%  '
$ws.Range('D43').Value = '1.532.53'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('E44').Value = '  +5.62%  '
$ws.Range('E45').Value = '  +7.16%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = "'0.0928"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = "'4.26"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.94%  '
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').Value = "'7.79"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').Value = '2.365.18'
$ws.Range('E50').Value = '  +3.29%  '
$ws.Range('E51').Value = '  +0.13%  '
